$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.704.07'
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").Value = '1.599.72'
$ws.Range("E3").Value = '  +0.42%  '

$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.43'
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("E6").Value = '  -0.33%  '

$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0618'
$ws.Range("E8").Value = '  +0.46%  '

$ws.Range("E9").Value = '  +1.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.51'
$ws.Range("E10").Value = '  +1.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0843'
$ws.Range("E11").Value = '  +0.83%  '

$ws.Range("D12").Value = '1.824.51'
$ws.Range("E12").Value = '  +0.41%  '

$ws.Range("D13").Value = '1.614.08'
$ws.Range("E13").Value = '  +0.54%  '

$ws.Range("E14").Value = '  +0.70%  '

$ws.Range("E15").Value = '  +0.79%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.34'
$ws.Range("E16").Value = '  +1.61%  '

$ws.Range("D17").Value = '26.684.64'
$ws.Range("E17").Value = '  +0.43%  '

$ws.Range("E18").Value = '  +3.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '209.41'
$ws.Range("E20").Value = '  +0.86%  '

$ws.Range("E21").Value = '  +3.72%  '

$ws.Range("E22").Value = '  +1.05%  '

$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.93'
$ws.Range("E24").Value = '  +1.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.97'
$ws.Range("E25").Value = '  -1.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("E26").Value = '  +0.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.12'
$ws.Range("E27").Value = '  +0.64%  '

$ws.Range("E28").Value = '  +0.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.34'
$ws.Range("E29").Value = '  +0.94%  '

$ws.Range("E30").Value = '  +2.75%  '

$ws.Range("E31").Value = '  +0.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.24'
$ws.Range("E32").Value = '  +0.97%  '

$ws.Range("E33").Value = '  +1.83%  '

$ws.Range("D34").Value = '1.289.22'
$ws.Range("E34").Value = '  +1.18%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.621'
$ws.Range("E35").Value = '  -4.68%  '

$ws.Range("E36").Value = '  +1.02%  '

$ws.Range("E37").Value = '  +0.74%  '

$ws.Range("E38").Value = '  +0.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.07'
$ws.Range("E39").Value = '  +16.40%  '

$ws.Range("E40").Value = '  -1.57%  '

$ws.Range("E41").Value = '  -0.74%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.785'
$ws.Range("E42").Value = '  +0.25%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.19'
$ws.Range("E43").Value = '  -0.36%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.09'
$ws.Range("E44").Value = '  -1.08%  '

$ws.Range("D45").Value = '1.737.03'
$ws.Range("E45").Value = '  +0.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.98'
$ws.Range("E46").Value = '  +1.87%  '

$ws.Range("E47").Value = '  -0.64%  '

$ws.Range("E48").Value = '  -0.67%  '

$ws.Range("E49").Value = '  +0.91%  '

$ws.Range("E51").Value = '  -0.04%  '
